# Auto-generated Excel COM-interop script to apply the Leviathan_Profits diff
# Updates computed price/profit columns (H-N) across sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17 (ALC)
$ws.Cells.Item(17, 8).Value = 69632.766  # H17: 53301.41 -> 69632.766
$ws.Cells.Item(17, 9).Value = 0  # I17: 99 -> 0
$ws.Cells.Item(17, 10).Value = 69632.766  # J17: 56626.562 -> 69632.766
$ws.Cells.Item(17, 11).Value = 0  # K17: 297 -> 0
$ws.Cells.Item(17, 12).Value = 208898.298  # L17: 169879.686 -> 208898.298
$ws.Cells.Item(17, 13).ClearContents()  # M17: was -129
$ws.Cells.Item(17, 14).Value = -209234.298  # N17: -170215.686 -> -209234.298

# Row 19 (ALC)
$ws.Cells.Item(19, 8).Value = 3994.5  # H19: 5094.6665 -> 3994.5
$ws.Cells.Item(19, 9).Value = 4891.5  # I19: 6226 -> 4891.5
$ws.Cells.Item(19, 10).Value = 3097.5  # J19: 3963.3333 -> 3097.5
$ws.Cells.Item(19, 11).Value = 4891.5  # K19: 6226 -> 4891.5
$ws.Cells.Item(19, 12).Value = 3097.5  # L19: 3963.3333 -> 3097.5
$ws.Cells.Item(19, 13).Value = -4716.5  # M19: -6051 -> -4716.5
$ws.Cells.Item(19, 14).Value = -3447.5  # N19: -4313.3333 -> -3447.5

# Row 51 (ALC)
$ws.Cells.Item(51, 8).Value = 6948811.5  # H51: 9263482 -> 6948811.5
$ws.Cells.Item(51, 9).Value = 4490.2  # I51: 3000 -> 4490.2
$ws.Cells.Item(51, 10).Value = 11909041  # J51: 10421042 -> 11909041
$ws.Cells.Item(51, 11).Value = 4490.2  # K51: 3000 -> 4490.2
$ws.Cells.Item(51, 12).Value = 11909041  # L51: 10421042 -> 11909041
$ws.Cells.Item(51, 13).Value = -4006.2  # M51: -2516 -> -4006.2
$ws.Cells.Item(51, 14).Value = -11910009  # N51: -10422010 -> -11910009

# Row 121 (ALC)
$ws.Cells.Item(121, 8).Value = 19998  # H121: 1500 -> 19998
$ws.Cells.Item(121, 10).Value = 19998  # J121: 1500 -> 19998
$ws.Cells.Item(121, 12).Value = 59994  # L121: 4500 -> 59994
$ws.Cells.Item(121, 14).Value = -63488  # N121: -7994 -> -63488

$ws = $wb.Worksheets.Item("ARM")
# Row 32 (ARM)
$ws.Cells.Item(32, 8).Value = 37410.586  # H32: 39413.13 -> 37410.586
$ws.Cells.Item(32, 9).Value = 22124.875  # I32: 23083.717 -> 22124.875
$ws.Cells.Item(32, 10).Value = 110782  # J32: 122874.555 -> 110782
$ws.Cells.Item(32, 11).Value = 22124.875  # K32: 23083.717 -> 22124.875
$ws.Cells.Item(32, 12).Value = 110782  # L32: 122874.555 -> 110782
$ws.Cells.Item(32, 13).Value = -21837.875  # M32: -22796.717 -> -21837.875
$ws.Cells.Item(32, 14).Value = -111356  # N32: -123448.555 -> -111356

# Row 45 (ARM)
$ws.Cells.Item(45, 8).Value = 845560.2  # H45: 845477.8 -> 845560.2
$ws.Cells.Item(45, 9).Value = 1687336.5  # I45: 1446575.9 -> 1687336.5
$ws.Cells.Item(45, 10).Value = 3783.8333  # J45: 3940.6 -> 3783.8333
$ws.Cells.Item(45, 11).Value = 1687336.5  # K45: 1446575.9 -> 1687336.5
$ws.Cells.Item(45, 12).Value = 3783.8333  # L45: 3940.6 -> 3783.8333
$ws.Cells.Item(45, 13).Value = -1686959.5  # M45: -1446198.9 -> -1686959.5
$ws.Cells.Item(45, 14).Value = -4537.8333  # N45: -4694.6 -> -4537.8333

# Row 61 (ARM)
$ws.Cells.Item(61, 8).Value = 2031.4  # H61: 1491.3125 -> 2031.4
$ws.Cells.Item(61, 9).Value = 1289.25  # I61: 1027.7693 -> 1289.25
$ws.Cells.Item(61, 10).Value = 5000  # J61: 3500 -> 5000
$ws.Cells.Item(61, 11).Value = 1289.25  # K61: 1027.7693 -> 1289.25
$ws.Cells.Item(61, 12).Value = 5000  # L61: 3500 -> 5000
$ws.Cells.Item(61, 13).Value = -1077.25  # M61: -815.7692999999999 -> -1077.25
$ws.Cells.Item(61, 14).Value = -5424  # N61: -3924 -> -5424

# Row 92 (ARM)
$ws.Cells.Item(92, 8).Value = 30550  # H92: 98000 -> 30550
$ws.Cells.Item(92, 10).Value = 30550  # J92: 98000 -> 30550
$ws.Cells.Item(92, 12).Value = 30550  # L92: 98000 -> 30550
$ws.Cells.Item(92, 14).Value = -35542  # N92: -102992 -> -35542

# Row 132 (ARM)
$ws.Cells.Item(132, 8).Value = 32293.285  # H132: 54888.25 -> 32293.285
$ws.Cells.Item(132, 9).Value = 42011  # I132: 101777.5 -> 42011
$ws.Cells.Item(132, 11).Value = 126033  # K132: 305332.5 -> 126033
$ws.Cells.Item(132, 13).Value = -123503  # M132: -302802.5 -> -123503

# Row 136 (ARM)
$ws.Cells.Item(136, 8).Value = 2031.4  # H136: 1491.3125 -> 2031.4
$ws.Cells.Item(136, 9).Value = 1289.25  # I136: 1027.7693 -> 1289.25
$ws.Cells.Item(136, 10).Value = 5000  # J136: 3500 -> 5000
$ws.Cells.Item(136, 11).Value = 3867.75  # K136: 3083.3079 -> 3867.75
$ws.Cells.Item(136, 12).Value = 15000  # L136: 10500 -> 15000
$ws.Cells.Item(136, 13).Value = -1317.75  # M136: -533.3078999999998 -> -1317.75
$ws.Cells.Item(136, 14).Value = -20100  # N136: -15600 -> -20100

$ws = $wb.Worksheets.Item("BSM")
# Row 22 (BSM)
$ws.Cells.Item(22, 8).Value = 529.3  # H22: 363.66666 -> 529.3
$ws.Cells.Item(22, 9).Value = 529.3  # I22: 399.75 -> 529.3
$ws.Cells.Item(22, 10).Value = 0  # J22: 75 -> 0
$ws.Cells.Item(22, 11).Value = 529.3  # K22: 399.75 -> 529.3
$ws.Cells.Item(22, 12).Value = 0  # L22: 75 -> 0
$ws.Cells.Item(22, 13).Value = -356.3  # M22: -226.75 -> -356.3
$ws.Cells.Item(22, 14).ClearContents()  # N22: was -421

# Row 107 (BSM)
$ws.Cells.Item(107, 8).Value = 11996.647  # H107: 10263.8 -> 11996.647
$ws.Cells.Item(107, 9).Value = 3371.4375  # I107: 2909.2632 -> 3371.4375
$ws.Cells.Item(107, 11).Value = 3371.4375  # K107: 2909.2632 -> 3371.4375
$ws.Cells.Item(107, 13).Value = -1451.4375  # M107: -989.2631999999999 -> -1451.4375

$ws = $wb.Worksheets.Item("CRP")
# Row 7 (CRP)
$ws.Cells.Item(7, 8).Value = 403.2  # H7: 390.80646 -> 403.2
$ws.Cells.Item(7, 9).Value = 256.6875  # I7: 242.70589 -> 256.6875
$ws.Cells.Item(7, 11).Value = 256.6875  # K7: 242.70589 -> 256.6875
$ws.Cells.Item(7, 13).Value = -143.6875  # M7: -129.70589 -> -143.6875

# Row 31 (CRP)
$ws.Cells.Item(31, 8).Value = 1599.234  # H31: 1576.3334 -> 1599.234
$ws.Cells.Item(31, 9).Value = 1268.1282  # I31: 1248.925 -> 1268.1282
$ws.Cells.Item(31, 11).Value = 1268.1282  # K31: 1248.925 -> 1268.1282
$ws.Cells.Item(31, 13).Value = -973.1282000000001  # M31: -953.925 -> -973.1282000000001

# Row 34 (CRP)
$ws.Cells.Item(34, 8).Value = 1599.234  # H34: 1576.3334 -> 1599.234
$ws.Cells.Item(34, 9).Value = 1268.1282  # I34: 1248.925 -> 1268.1282
$ws.Cells.Item(34, 11).Value = 1268.1282  # K34: 1248.925 -> 1268.1282
$ws.Cells.Item(34, 13).Value = -1066.1282  # M34: -1046.925 -> -1066.1282

# Row 58 (CRP)
$ws.Cells.Item(58, 8).Value = 10844.363  # H58: 10935.637 -> 10844.363
$ws.Cells.Item(58, 9).Value = 34166.668  # I58: 50750 -> 34166.668
$ws.Cells.Item(58, 10).Value = 2098.5  # J58: 2088 -> 2098.5
$ws.Cells.Item(58, 11).Value = 34166.668  # K58: 50750 -> 34166.668
$ws.Cells.Item(58, 12).Value = 2098.5  # L58: 2088 -> 2098.5
$ws.Cells.Item(58, 13).Value = -33963.668  # M58: -50547 -> -33963.668
$ws.Cells.Item(58, 14).Value = -2504.5  # N58: -2494 -> -2504.5

# Row 132 (CRP)
$ws.Cells.Item(132, 8).Value = 3126.0908  # H132: 3094.739 -> 3126.0908
$ws.Cells.Item(132, 10).Value = 4100.5  # J132: 3761.4 -> 4100.5
$ws.Cells.Item(132, 12).Value = 12301.5  # L132: 11284.2 -> 12301.5
$ws.Cells.Item(132, 14).Value = -17361.5  # N132: -16344.2 -> -17361.5

# Row 134 (CRP)
$ws.Cells.Item(134, 8).Value = 4197.273  # H134: 4117.769 -> 4197.273
$ws.Cells.Item(134, 9).Value = 4020  # I134: 3958.2727 -> 4020
$ws.Cells.Item(134, 11).Value = 12060  # K134: 11874.8181 -> 12060
$ws.Cells.Item(134, 13).Value = -9525  # M134: -9339.8181 -> -9525

# Row 136 (CRP)
$ws.Cells.Item(136, 8).Value = 10844.363  # H136: 10935.637 -> 10844.363
$ws.Cells.Item(136, 9).Value = 34166.668  # I136: 50750 -> 34166.668
$ws.Cells.Item(136, 10).Value = 2098.5  # J136: 2088 -> 2098.5
$ws.Cells.Item(136, 11).Value = 102500.004  # K136: 152250 -> 102500.004
$ws.Cells.Item(136, 12).Value = 6295.5  # L136: 6264 -> 6295.5
$ws.Cells.Item(136, 13).Value = -99950.00399999999  # M136: -149700 -> -99950.00399999999
$ws.Cells.Item(136, 14).Value = -11395.5  # N136: -11364 -> -11395.5

$ws = $wb.Worksheets.Item("CUL")
# Row 37 (CUL)
$ws.Cells.Item(37, 8).Value = 333400000  # H37: 1000000000 -> 333400000
$ws.Cells.Item(37, 10).Value = 333400000  # J37: 1000000000 -> 333400000
$ws.Cells.Item(37, 12).Value = 1000200000  # L37: 3000000000 -> 1000200000
$ws.Cells.Item(37, 14).Value = -1000200224  # N37: -3000000224 -> -1000200224

# Row 138 (CUL)
$ws.Cells.Item(138, 8).Value = 3070.8125  # H138: 3429.5386 -> 3070.8125
$ws.Cells.Item(138, 9).Value = 1158.1666  # I138: 1224.9 -> 1158.1666
$ws.Cells.Item(138, 10).Value = 8808.75  # J138: 10778.333 -> 8808.75
$ws.Cells.Item(138, 11).Value = 3474.4998  # K138: 3674.7 -> 3474.4998
$ws.Cells.Item(138, 12).Value = 26426.25  # L138: 32334.999 -> 26426.25
$ws.Cells.Item(138, 13).Value = 1665.5002  # M138: 1465.3 -> 1665.5002
$ws.Cells.Item(138, 14).Value = -36706.25  # N138: -42614.999 -> -36706.25

$ws = $wb.Worksheets.Item("GSM")
# Row 122 (GSM)
$ws.Cells.Item(122, 8).Value = 1313.8125  # H122: 1101.2727 -> 1313.8125
$ws.Cells.Item(122, 9).Value = 1428.9166  # I122: 1177.1428 -> 1428.9166
$ws.Cells.Item(122, 11).Value = 4286.7498  # K122: 3531.4284 -> 4286.7498
$ws.Cells.Item(122, 13).Value = -1836.7498  # M122: -1081.4284 -> -1836.7498

# Row 132 (GSM)
$ws.Cells.Item(132, 8).Value = 2900  # H132: 1885.25 -> 2900
$ws.Cells.Item(132, 9).Value = 900  # I132: 885 -> 900
$ws.Cells.Item(132, 10).Value = 4900  # J132: 2885.5 -> 4900
$ws.Cells.Item(132, 11).Value = 2700  # K132: 2655 -> 2700
$ws.Cells.Item(132, 12).Value = 14700  # L132: 8656.5 -> 14700
$ws.Cells.Item(132, 13).Value = -170  # M132: -125 -> -170
$ws.Cells.Item(132, 14).Value = -19760  # N132: -13716.5 -> -19760

$ws = $wb.Worksheets.Item("LTW")
# Row 55 (LTW)
$ws.Cells.Item(55, 8).Value = 486.13333  # H55: 510.18182 -> 486.13333
$ws.Cells.Item(55, 9).Value = 502  # I55: 476.75 -> 502
$ws.Cells.Item(55, 10).Value = 472.25  # J55: 599.3333 -> 472.25
$ws.Cells.Item(55, 11).Value = 502  # K55: 476.75 -> 502
$ws.Cells.Item(55, 12).Value = 472.25  # L55: 599.3333 -> 472.25
$ws.Cells.Item(55, 13).Value = -329  # M55: -303.75 -> -329
$ws.Cells.Item(55, 14).Value = -818.25  # N55: -945.3333 -> -818.25

# Row 61 (LTW)
$ws.Cells.Item(61, 8).Value = 67988.8  # H61: 56990.668 -> 67988.8
$ws.Cells.Item(61, 9).Value = 1652.75  # I61: 1722.2 -> 1652.75
$ws.Cells.Item(61, 11).Value = 1652.75  # K61: 1722.2 -> 1652.75
$ws.Cells.Item(61, 13).Value = -1450.75  # M61: -1520.2 -> -1450.75

# Row 107 (LTW)
$ws.Cells.Item(107, 8).Value = 2571.2856  # H107: 3199.6667 -> 2571.2856
$ws.Cells.Item(107, 9).Value = 2571.2856  # I107: 3199.6667 -> 2571.2856
$ws.Cells.Item(107, 11).Value = 2571.2856  # K107: 3199.6667 -> 2571.2856
$ws.Cells.Item(107, 13).Value = -651.2856000000002  # M107: -1279.6667 -> -651.2856000000002

# Row 113 (LTW)
$ws.Cells.Item(113, 8).Value = 67988.8  # H113: 56990.668 -> 67988.8
$ws.Cells.Item(113, 9).Value = 1652.75  # I113: 1722.2 -> 1652.75
$ws.Cells.Item(113, 11).Value = 1652.75  # K113: 1722.2 -> 1652.75
$ws.Cells.Item(113, 13).Value = 517.25  # M113: 447.8 -> 517.25

$ws = $wb.Worksheets.Item("WVR")
# Row 111 (WVR)
$ws.Cells.Item(111, 8).Value = 21744  # H111: 11248 -> 21744
$ws.Cells.Item(111, 10).Value = 21744  # J111: 11248 -> 21744
$ws.Cells.Item(111, 12).Value = 21744  # L111: 11248 -> 21744
$ws.Cells.Item(111, 14).Value = -29924  # N111: -19428 -> -29924

# Row 132 (WVR)
$ws.Cells.Item(132, 8).Value = 17585.117  # H132: 17626.295 -> 17585.117
$ws.Cells.Item(132, 9).Value = 17217.455  # I132: 17281.092 -> 17217.455
$ws.Cells.Item(132, 11).Value = 51652.36500000001  # K132: 51843.276 -> 51652.36500000001
$ws.Cells.Item(132, 13).Value = -49122.36500000001  # M132: -49313.276 -> -49122.36500000001

# Row 136 (WVR)
$ws.Cells.Item(136, 8).Value = 751.2857  # H136: 911.8 -> 751.2857
$ws.Cells.Item(136, 9).Value = 809.8333  # I136: 911.8 -> 809.8333
$ws.Cells.Item(136, 10).Value = 400  # J136: 0 -> 400
$ws.Cells.Item(136, 11).Value = 2429.4999  # K136: 2735.4 -> 2429.4999
$ws.Cells.Item(136, 12).Value = 1200  # L136: 0 -> 1200
$ws.Cells.Item(136, 13).Value = 120.5001000000002  # M136: -185.3999999999996 -> 120.5001000000002
$ws.Cells.Item(136, 14).Value = -6300  # N136: None -> -6300
